# Auto-generated edit script applying scheduled market-data refresh to the
# per-job-sheet Leve profit tables (columns H:N) as captured in the commit diff.
$wb = $excel.ActiveWorkbook

# --- Sheet "ALC" ---
$ws = $wb.Worksheets.Item("ALC")

# Row 17
$ws.Range("H17").Value = 1862.5454
$ws.Range("I17").Value = 1453.3334
$ws.Range("J17").Value = 2353.6
$ws.Range("K17").Value = 4360.0002
$ws.Range("L17").Value = 7060.799999999999
$ws.Range("M17").Value = -4192.0002
$ws.Range("N17").Value = -7396.799999999999

# Row 32
$ws.Range("H32").Value = 12453.55
$ws.Range("I32").Value = 10934.363
$ws.Range("K32").Value = 10934.363
$ws.Range("M32").Value = -10608.363

# Row 40
$ws.Range("H40").Value = 1815.8334
$ws.Range("I40").Value = 1839
$ws.Range("J40").Value = 1700
$ws.Range("K40").Value = 1839
$ws.Range("L40").Value = 1700
$ws.Range("M40").Value = -1664
$ws.Range("N40").Value = -2050

# Row 86
$ws.Range("H86").Value = 2618.5
$ws.Range("J86").Value = 2766.9092
$ws.Range("L86").Value = 2766.9092
$ws.Range("N86").Value = -5012.9092

# Row 89
$ws.Range("H89").Value = 2618.5
$ws.Range("J89").Value = 2766.9092
$ws.Range("L89").Value = 13834.546
$ws.Range("N89").Value = -25066.546

# Row 100
$ws.Range("H100").Value = 80540.125
$ws.Range("I100").Value = 151685.5
$ws.Range("K100").Value = 151685.5
$ws.Range("M100").Value = -151144.5

# Row 113
$ws.Range("H113").Value = 4265.5835
$ws.Range("I113").Value = 3098
$ws.Range("K113").Value = 3098
$ws.Range("M113").Value = 156

# Row 135
$ws.Range("H135").Value = 15561.237
$ws.Range("I135").Value = 1637.7307
$ws.Range("J135").Value = 45728.832
$ws.Range("K135").Value = 14739.5763
$ws.Range("L135").Value = 411559.488
$ws.Range("M135").Value = -12204.5763
$ws.Range("N135").Value = -416629.488

# Row 138
$ws.Range("H138").Value = 2957.84
$ws.Range("I138").Value = 2008.4117
$ws.Range("J138").Value = 4975.375
$ws.Range("K138").Value = 6025.2351
$ws.Range("L138").Value = 14926.125
$ws.Range("M138").Value = -885.2350999999999
$ws.Range("N138").Value = -25206.125

# --- Sheet "ARM" ---
$ws = $wb.Worksheets.Item("ARM")

# Row 32
$ws.Range("H32").Value = 18610.941
$ws.Range("I32").Value = 19080.984
$ws.Range("K32").Value = 19080.984
$ws.Range("M32").Value = -18793.984

# Row 45
$ws.Range("H45").Value = 3361.9714
$ws.Range("I45").Value = 2310.4546
$ws.Range("K45").Value = 2310.4546
$ws.Range("M45").Value = -1933.4546

# Row 61
$ws.Range("H61").Value = 2906.4
$ws.Range("I61").Value = 2009.5714
$ws.Range("K61").Value = 2009.5714
$ws.Range("M61").Value = -1797.5714

# Row 119
$ws.Range("H119").Value = 0
$ws.Range("I119").Value = 0
$ws.Range("K119").Value = 0
$ws.Range("M119").ClearContents()

# Row 124
$ws.Range("H124").Value = 50976
$ws.Range("J124").Value = 50976
$ws.Range("L124").Value = 50976
$ws.Range("N124").Value = -60796

# Row 125
$ws.Range("H125").Value = 34410
$ws.Range("J125").Value = 34410
$ws.Range("L125").Value = 34410
$ws.Range("N125").Value = -44250

# Row 136
$ws.Range("H136").Value = 2906.4
$ws.Range("I136").Value = 2009.5714
$ws.Range("K136").Value = 6028.7142
$ws.Range("M136").Value = -3478.7142

# Row 141
$ws.Range("H141").Value = 75000
$ws.Range("J141").Value = 75000
$ws.Range("L141").Value = 75000
$ws.Range("N141").Value = -85360

# --- Sheet "BSM" ---
$ws = $wb.Worksheets.Item("BSM")

# Row 20
$ws.Range("H20").Value = 15595.177
$ws.Range("I20").Value = 21711.25
$ws.Range("J20").Value = 916.6
$ws.Range("K20").Value = 21711.25
$ws.Range("L20").Value = 916.6
$ws.Range("M20").Value = -21464.25
$ws.Range("N20").Value = -1410.6

# Row 22
$ws.Range("H22").Value = 243.21428
$ws.Range("I22").Value = 250.63637
$ws.Range("K22").Value = 250.63637
$ws.Range("M22").Value = -77.63637

# Row 94
$ws.Range("H94").Value = 949.7143
$ws.Range("I94").Value = 975
$ws.Range("K94").Value = 975
$ws.Range("M94").Value = -524

# Row 107
$ws.Range("H107").Value = 29377.555
$ws.Range("I107").Value = 43301.832
$ws.Range("J107").Value = 1529
$ws.Range("K107").Value = 43301.832
$ws.Range("L107").Value = 1529
$ws.Range("M107").Value = -41381.832
$ws.Range("N107").Value = -5369

# --- Sheet "CRP" ---
$ws = $wb.Worksheets.Item("CRP")

# Row 31
$ws.Range("H31").Value = 2943952.8
$ws.Range("I31").Value = 3706305.5
$ws.Range("J31").Value = 3448.5715
$ws.Range("K31").Value = 3706305.5
$ws.Range("L31").Value = 3448.5715
$ws.Range("M31").Value = -3706010.5
$ws.Range("N31").Value = -4038.5715

# Row 34
$ws.Range("H34").Value = 2943952.8
$ws.Range("I34").Value = 3706305.5
$ws.Range("J34").Value = 3448.5715
$ws.Range("K34").Value = 3706305.5
$ws.Range("L34").Value = 3448.5715
$ws.Range("M34").Value = -3706103.5
$ws.Range("N34").Value = -3852.5715

# Row 92
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

# Row 99
$ws.Range("H99").Value = 7473.3335
$ws.Range("I99").Value = 5783.125
$ws.Range("K99").Value = 5783.125
$ws.Range("M99").Value = -4285.125

# Row 107
$ws.Range("H107").Value = 949.9231
$ws.Range("I107").Value = 773.1429000000001
$ws.Range("J107").Value = 1156.1666
$ws.Range("K107").Value = 773.1429000000001
$ws.Range("L107").Value = 1156.1666
$ws.Range("M107").Value = 1146.8571
$ws.Range("N107").Value = -4996.1666

# Row 126
$ws.Range("H126").Value = 7473.3335
$ws.Range("I126").Value = 5783.125
$ws.Range("K126").Value = 17349.375
$ws.Range("M126").Value = -14879.375

# --- Sheet "CUL" ---
$ws = $wb.Worksheets.Item("CUL")

# Row 51
$ws.Range("H51").Value = 3648.1428
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 3648.1428
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 10944.4284
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -11864.4284

# Row 80
$ws.Range("H80").Value = 4000
$ws.Range("I80").Value = 4000
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 12000
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -11064
$ws.Range("N80").ClearContents()

# Row 83
$ws.Range("H83").Value = 4000
$ws.Range("I83").Value = 4000
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 36000
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -31320
$ws.Range("N83").ClearContents()

# Row 87
$ws.Range("H87").Value = 7008.8184
$ws.Range("I87").Value = 2637.125
$ws.Range("K87").Value = 7911.375
$ws.Range("M87").Value = -6663.375

# Row 90
$ws.Range("H90").Value = 7008.8184
$ws.Range("I90").Value = 2637.125
$ws.Range("K90").Value = 23734.125
$ws.Range("M90").Value = -17494.125

# Row 103
$ws.Range("H103").Value = 185.5
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()

# Row 113
$ws.Range("H113").Value = 1421.3529
$ws.Range("J113").Value = 1697.8334
$ws.Range("L113").Value = 5093.5002
$ws.Range("N113").Value = -9433.5002

# Row 129
$ws.Range("H129").Value = 4060.5
$ws.Range("I129").Value = 1587.7142
$ws.Range("J129").Value = 6533.2856
$ws.Range("K129").Value = 4763.142599999999
$ws.Range("L129").Value = 19599.8568
$ws.Range("M129").Value = 236.8574000000008
$ws.Range("N129").Value = -29599.8568

# Row 136
$ws.Range("H136").Value = 3632.75
$ws.Range("I136").Value = 2332
$ws.Range("K136").Value = 6996
$ws.Range("M136").Value = -1896

# Row 138
$ws.Range("H138").Value = 27555
$ws.Range("I138").Value = 51110.5
$ws.Range("K138").Value = 153331.5
$ws.Range("M138").Value = -148191.5

# --- Sheet "GSM" ---
$ws = $wb.Worksheets.Item("GSM")

# Row 80
$ws.Range("H80").Value = 7963.2
$ws.Range("I80").Value = 4599.3335
$ws.Range("K80").Value = 4599.3335
$ws.Range("M80").Value = -3601.3335

# Row 83
$ws.Range("H83").Value = 7963.2
$ws.Range("I83").Value = 4599.3335
$ws.Range("K83").Value = 22996.6675
$ws.Range("M83").Value = -18004.6675

# Row 107
$ws.Range("H107").Value = 303.5
$ws.Range("I107").Value = 217
$ws.Range("K107").Value = 217
$ws.Range("M107").Value = 1703

# --- Sheet "LTW" ---
$ws = $wb.Worksheets.Item("LTW")

# Row 21
$ws.Range("H21").Value = 16171.25
$ws.Range("J21").Value = 16561.666
$ws.Range("L21").Value = 16561.666
$ws.Range("N21").Value = -16909.666

# Row 22
$ws.Range("H22").Value = 1222
$ws.Range("I22").Value = 921.1111
$ws.Range("J22").Value = 1447.6666
$ws.Range("K22").Value = 921.1111
$ws.Range("L22").Value = 1447.6666
$ws.Range("M22").Value = -626.1111
$ws.Range("N22").Value = -2037.6666

# Row 27
$ws.Range("H27").Value = 1222
$ws.Range("I27").Value = 921.1111
$ws.Range("J27").Value = 1447.6666
$ws.Range("K27").Value = 921.1111
$ws.Range("L27").Value = 1447.6666
$ws.Range("M27").Value = -814.1111
$ws.Range("N27").Value = -1661.6666

# Row 40
$ws.Range("H40").Value = 3130.2222
$ws.Range("I40").Value = 3440.7144
$ws.Range("K40").Value = 3440.7144
$ws.Range("M40").Value = -3304.7144

# Row 74
$ws.Range("H74").Value = 50398
$ws.Range("J74").Value = 57333.332
$ws.Range("L74").Value = 57333.332
$ws.Range("N74").Value = -59329.332

# Row 77
$ws.Range("H77").Value = 50398
$ws.Range("J77").Value = 57333.332
$ws.Range("L77").Value = 171999.996
$ws.Range("N77").Value = -181983.996

# Row 122
$ws.Range("H122").Value = 7148.913
$ws.Range("I122").Value = 7119.909
$ws.Range("K122").Value = 21359.727
$ws.Range("M122").Value = -18909.727

# Row 136
$ws.Range("H136").Value = 7206.6
$ws.Range("I136").Value = 7925
$ws.Range("J136").Value = 6727.6665
$ws.Range("K136").Value = 23775
$ws.Range("L136").Value = 20182.9995
$ws.Range("M136").Value = -21225
$ws.Range("N136").Value = -25282.9995

# --- Sheet "WVR" ---
$ws = $wb.Worksheets.Item("WVR")

# Row 136
$ws.Range("H136").Value = 16983.018
$ws.Range("I136").Value = 20323.88
$ws.Range("J136").Value = 5290
$ws.Range("K136").Value = 60971.64
$ws.Range("L136").Value = 15870
$ws.Range("M136").Value = -58421.64
$ws.Range("N136").Value = -20970
